# faturamento_diario.xlsx — "atualizei dados bibi e add"
#
# The June/May/April "Dia" series previously started counting at 1 right
# away under July, i.e. July (month 7) was missing its 31st day. That
# missing day (31/07/2025, total_venda 27424.66) needs to be inserted,
# pushing the June/May data down by one row. At the same time, the whole
# April 2025 (month 4) block at the bottom of the sheet is removed,
# because it's no longer part of the tracked period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the missing July 31st row right before the old row 32 (which
#    used to be "day 1" of June) — this shifts June/May down by one row.
$ws.Rows(32).Insert()
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 27424.66
$ws.Range("C32").Value = 7
$ws.Range("D32").Value = 2025
$ws.Range("E32").Value = "07/2025"

# 2) Drop the April 2025 block entirely. Before the insert above it sat
#    at rows 93:122 (30 rows); after shifting everything down by one row
#    it now lives at 94:123.
$ws.Range("A94:E123").EntireRow.Delete()
